$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 32.21267
$ws.Range("H2").Value = 96.63801000000001
$ws.Range("I2").Value = 0.7096649552378644
$ws.Range("J2").Value = 0.7096649552378644
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 4693.410006701594
$ws.Range("R2").Value = 42240.69006031435
$ws.Range("S2").Value = 0.2033855585143369
$ws.Range("T2").Value = 0.2033855585143369
$ws.Range("G3").Value = 32.21267
$ws.Range("H3").Value = 96.63801000000001
$ws.Range("I3").Value = 0.7096649552378644
$ws.Range("J3").Value = 0.7096649552378644
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 5437.489118099455
$ws.Range("R3").Value = 48937.40206289509
$ws.Range("S3").Value = 0.2356296934683294
$ws.Range("T3").Value = 0.2356296934683294
$ws.Range("G4").Value = 32.21267
$ws.Range("H4").Value = 96.63801000000001
$ws.Range("I4").Value = 0.7096649552378644
$ws.Range("J4").Value = 0.7096649552378644
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 4127.285538646294
$ws.Range("R4").Value = 37145.56984781665
$ws.Range("S4").Value = 0.1788529604758848
$ws.Range("T4").Value = 0.1788529604758848
$ws.Range("G5").Value = 32.21267
$ws.Range("H5").Value = 96.63801000000001
$ws.Range("I5").Value = 0.7096649552378644
$ws.Range("J5").Value = 0.7096649552378644
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 2118.33993667093
$ws.Range("R5").Value = 19065.05943003837
$ws.Range("S5").Value = 0.09179674277931346
$ws.Range("T5").Value = 0.09179674277931346
$ws.Range("I6").Value = 0.2527239295880077
$ws.Range("J6").Value = 0.2527239295880077
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 1671.404246900901
$ws.Range("R6").Value = 15042.63822210811
$ws.Range("S6").Value = 0.07242910501613622
$ws.Range("T6").Value = 0.07242910501613622
$ws.Range("I7").Value = 0.2527239295880077
$ws.Range("J7").Value = 0.2527239295880077
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.08391179756223736
$ws.Range("T7").Value = 0.08391179756223736
$ws.Range("I8").Value = 0.2527239295880077
$ws.Range("J8").Value = 0.2527239295880077
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 1469.797560327374
$ws.Range("R8").Value = 13228.17804294637
$ws.Range("S8").Value = 0.06369262376041102
$ws.Range("T8").Value = 0.06369262376041102
$ws.Range("I9").Value = 0.2527239295880077
$ws.Range("J9").Value = 0.2527239295880077
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 754.3773847748326
$ws.Range("R9").Value = 6789.396462973494
$ws.Range("S9").Value = 0.03269040324922313
$ws.Range("T9").Value = 0.03269040324922313
$ws.Range("G10").Value = 1.279382333333333
$ws.Range("H10").Value = 3.838147
$ws.Range("I10").Value = 0.02818558059040478
$ws.Range("J10").Value = 0.02818558059040478
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 186.4069586800442
$ws.Range("R10").Value = 1677.662628120398
$ws.Range("S10").Value = 0.008077811942269159
$ws.Range("T10").Value = 0.008077811942269159
$ws.Range("G11").Value = 1.279382333333333
$ws.Range("H11").Value = 3.838147
$ws.Range("I11").Value = 0.02818558059040478
$ws.Range("J11").Value = 0.02818558059040478
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 215.9593574636529
$ws.Range("R11").Value = 1943.634217172876
$ws.Range("S11").Value = 0.009358443961091376
$ws.Range("T11").Value = 0.009358443961091376
$ws.Range("G12").Value = 1.279382333333333
$ws.Range("H12").Value = 3.838147
$ws.Range("I12").Value = 0.02818558059040478
$ws.Range("J12").Value = 0.02818558059040478
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 163.9223387184676
$ws.Range("R12").Value = 1475.301048466208
$ws.Range("S12").Value = 0.007103457052681813
$ws.Range("T12").Value = 0.007103457052681813
$ws.Range("G13").Value = 1.279382333333333
$ws.Range("H13").Value = 3.838147
$ws.Range("I13").Value = 0.02818558059040478
$ws.Range("J13").Value = 0.02818558059040478
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 84.13356269353766
$ws.Range("R13").Value = 757.2020642418389
$ws.Range("S13").Value = 0.003645867634362437
$ws.Range("T13").Value = 0.003645867634362437
$ws.Range("G14").Value = 0.4278379999999999
$ws.Range("H14").Value = 1.283514
$ws.Range("I14").Value = 0.009425534583723031
$ws.Range("J14").Value = 0.009425534583723031
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 62.33631519669733
$ws.Range("R14").Value = 561.026836770276
$ws.Range("S14").Value = 0.002701299537841999
$ws.Range("T14").Value = 0.002701299537841999
$ws.Range("G15").Value = 0.4278379999999999
$ws.Range("H15").Value = 1.283514
$ws.Range("I15").Value = 0.009425534583723031
$ws.Range("J15").Value = 0.009425534583723031
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 72.21892718950133
$ws.Range("R15").Value = 649.970344705512
$ws.Range("S15").Value = 0.003129555444925959
$ws.Range("T15").Value = 0.003129555444925959
$ws.Range("G16").Value = 0.4278379999999999
$ws.Range("H16").Value = 1.283514
$ws.Range("I16").Value = 0.009425534583723031
$ws.Range("J16").Value = 0.009425534583723031
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 54.81723776027733
$ws.Range("R16").Value = 493.3551398424959
$ws.Range("S16").Value = 0.002375465706632874
$ws.Range("T16").Value = 0.002375465706632874
$ws.Range("G17").Value = 0.4278379999999999
$ws.Range("H17").Value = 1.283514
$ws.Range("I17").Value = 0.009425534583723031
$ws.Range("J17").Value = 0.009425534583723031
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 28.13508851720199
$ws.Range("R17").Value = 253.215796654818
$ws.Range("S17").Value = 0.0012192138943222
$ws.Range("T17").Value = 0.0012192138943222
